$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.271.80"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "3.514.45"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'586.34"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").Value = "'133.17"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").Value = "3.517.59"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.490"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").Value = "'0.124"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "'0.383"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "4.101.10"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "'27.74"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "'0.0000180"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "3.508.20"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "64.154.10"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "'10.19"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("D20").Value = "'14.45"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "'5.71"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'386.01"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "'0.580"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "3.653.30"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'73.35"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'0.0000116"
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("D28").Value = "'1.58"
$ws.Range("E28").Value = "  -2.53%  "
$ws.Range("D29").Value = "'7.59"
$ws.Range("E29").Value = "  -2.52%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "'2.27"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "'8.35"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").Value = "3.518.91"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'23.84"
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'5.40"
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("D38").Value = "'1.60"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").Value = "'6.99"
$ws.Range("E39").Value = "  +1.96%  "
$ws.Range("D40").Value = "'161.54"
$ws.Range("E40").Value = "  -4.84%  "
$ws.Range("D41").Value = "'0.0806"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("D42").Value = "'0.816"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").Value = "'26.13"
$ws.Range("E43").Value = "  +3.21%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").Value = "'1.24"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.998"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").Value = "'41.68"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").Value = "'4.43"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").Value = "'1.66"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "'6.89"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").Value = "2.435.02"
$ws.Range("E50").Value = "  +2.73%  "
$ws.Range("D51").Value = "'0.0269"
$ws.Range("E51").Value = "  +1.05%  "
